$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Copy-RowFormat($srcRow, $dstRow) {
    $ws.Range("A$srcRow`:F$srcRow").Copy() | Out-Null
    $ws.Range("A$dstRow`:F$dstRow").PasteSpecial(-4122) | Out-Null
}

function Copy-RowValues($srcRow, $dstRow) {
    for ($col = 1; $col -le 6; $col++) {
        $srcCell = $ws.Cells.Item($srcRow, $col)
        $dstCell = $ws.Cells.Item($dstRow, $col)
        $dstCell.Value2 = $srcCell.Value2
    }
}

# ------------------------------------------------------------------
# 1) Make room for the new "Tasa de apertura" indicator row: push the
#    "Construcción" block (old rows 74-76) one row down (new rows
#    75-77), shifting bottom-up so sources are never overwritten
#    before they are read. This keeps the style table untouched
#    (no Insert() call, which would otherwise mint new style ids).
# ------------------------------------------------------------------
Copy-RowFormat 76 77
Copy-RowValues 76 77

Copy-RowFormat 75 76
Copy-RowValues 75 76

Copy-RowFormat 74 75
Copy-RowValues 74 75

# ------------------------------------------------------------------
# 2) Populate the freed-up row 74 with the new dataset, reusing the
#    "highlighted" (new-data) style that rows 6/58/59 used to have.
# ------------------------------------------------------------------
Copy-RowFormat 58 74
$ws.Range("A74").Value2 = "Tasa de apertura"
$ws.Range("B74").Value2 = "%"
$ws.Range("C74").Value2 = ""
$ws.Range("D74").Value2 = "Fuente: ICANE a partir de Estadisitca de Comercio Exterior de la AEAT, Contabilidad Nacional Trimestral de España del INE y Contabilidad Trimestral de Cantabria Base 2015 del ICANE"
$ws.Range("E74").Value2 = "Tasa de apertura= (saldo comercial/PIB)*100. Para una mejor interpretación la tasa de variación se da en términos absolutos al tratarse de un saldo de porcentajes. Datos provisionales"
$ws.Range("F74").Value2 = ""

# ------------------------------------------------------------------
# 3) Normalize the formatting ("red flag" style) of rows 6, 58 and 59
#    back to the regular black style used elsewhere in the sheet
#    (style 5/6 instead of 13/10). Copy formats only from neighboring
#    rows that already use the desired style, so no new style
#    entries get created.
# ------------------------------------------------------------------
Copy-RowFormat 7 6

Copy-RowFormat 57 58
Copy-RowFormat 57 59

$ws.Application.CutCopyMode = $false

# ------------------------------------------------------------------
# 4) Restore the sheet view the same way the workbook was left after
#    the edit (no frozen/top-left scroll position, new selection).
# ------------------------------------------------------------------
$ws.Application.Goto($ws.Range("A1"), $true)
$ws.Range("E81").Select() | Out-Null
